$d = $word.ActiveDocument
$pos = $d.Content.End
$r = $d.Range($pos, $pos)
$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:lastRenderedPageBreak/><w:t>Incorrect Multiplication: Remember that in multiples of 2, simply double the number.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Incorrect Multiplication: The trick with multiples of 3 is to </w:t></w:r><w:r><w:t>double the number</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> and then add the original number.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Incorrect Multiplication: In multiples of 4, double the number, and then double it again.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Incorrect Multiplication: </w:t></w:r><w:r><w:t>If you are having trouble with multiples of 5, try multiplying the number by 10, and then half it.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Incorrect Multiplication: A good way to solve multiples of 6 is to multiply the number by 5, and then add the original number.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Incorrect Multiplication: A good way to solve multiples of 7 is to multiply the number by 5, and then add the original number twice.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Incorrect Multiplication: If you are having trouble with multiples of 8, try multiplying the number by 2 three times.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Incorrect Multiplication: One way to do multiples of 9 is to multiply the number by 10, and then subtract it by the original number.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Incorrect number of zeroes at the end! </w:t></w:r><w:r><w:t>Looks like you missed a few zeroes.</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">Incorrect number of zeroes at the end! Looks like you </w:t></w:r><w:r><w:t>put in too</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>many</w:t></w:r><w:r><w:t xml:space="preserve"> zeroes.</w:t></w:r></w:p>')
